# Generate Report for Handoff
#
# Refresh the "latest handoff" timestamps for the da596fb6-... row (row 7)
# across the Overview, zh-cn and de-de sheets, simulating a freshly
# generated handoff report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$wsOverview.Range("G7").Value = "2016-08-31 18:48:40"

# zh-cn sheet: column H = "Latest Handoff Datetime"
$wsZhCn.Range("H7").Value = "2016-08-31 18:48:35"

# de-de sheet: column H = "Latest Handoff Datetime"
$wsDeDe.Range("H7").Value = "2016-08-31 18:48:40"
